$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The new "QPSK 600" preset block (rows 56-60) is a structural clone of the
# existing QPSK 600 block at rows 20-24, but with a narrower margin (25%
# instead of 75%) and an extra "adjust limits" note. Grab the borders /
# number-formats from the template block first, cell-range by cell-range so
# we only touch the cells that are actually populated in the template (and
# don't leave stray blank styled cells behind).
$ws.Range("A20:J20").Copy() | Out-Null
$ws.Range("A56").PasteSpecial(-4122) | Out-Null
$ws.Range("A21:J21").Copy() | Out-Null
$ws.Range("A57").PasteSpecial(-4122) | Out-Null
$ws.Range("A22:E22").Copy() | Out-Null
$ws.Range("A58").PasteSpecial(-4122) | Out-Null
$ws.Range("J22").Copy() | Out-Null
$ws.Range("J58").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Copy() | Out-Null
$ws.Range("A59").PasteSpecial(-4122) | Out-Null
$ws.Range("B23").Copy() | Out-Null
$ws.Range("B59").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Copy() | Out-Null
$ws.Range("D59").PasteSpecial(-4122) | Out-Null
$ws.Range("J23").Copy() | Out-Null
$ws.Range("J59").PasteSpecial(-4122) | Out-Null
$ws.Range("A24:J24").Copy() | Out-Null
$ws.Range("A60").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 55 is the blank "thick bottom border" divider row that separates this
# new block from the one above it (same pattern as rows 7/13/19/25/31/37/43/49).
$ws.Rows.Item(55).RowHeight = 15.75

# Row 56 - "NCO Design Samp Rate" header line
$ws.Range("A56").Formula = "NCO Design Samp Rate"
$ws.Range("B56").Formula = 65536

# Row 57 - "Sample Rate" line
$ws.Range("A57").Formula = "Sample Rate"
$ws.Range("B57").Formula = 14400

# Row 58 - "I scale bits" / "p scale bits" / "max integral product bits" line
$ws.Range("B58").Formula = "I scale bits"
$ws.Range("C58").Formula = 20
$ws.Range("D58").Formula = "p scale bits"
$ws.Range("E58").Formula = 15
$ws.Range("J58").Formula = "max integral product bits"
$ws.Range("L58").Formula = "adjust for 13 hz offset"

# Row 59 - "integral gain" / "p gain" line
$ws.Range("B59").Formula = "integral gain"
$ws.Range("D59").Formula = "p gain"
$ws.Range("L59").Formula = "adjust limits"

# Row 60 - "QPSK 600" result line, with the adjusted (25%) margin
$ws.Range("A60").Formula = "QPSK 600"
$ws.Range("B60").Formula = 0.000031
$ws.Range("C60").Formula = '=ROUND(POWER(2,$C$16)*B60, 0)'
$ws.Range("D60").Formula = 0.02
$ws.Range("E60").Formula = '=ROUND(POWER(2,$E$4)*D60, 0)'
$ws.Range("F60").Formula = 13
$ws.Range("G60").Formula = 0.25
$ws.Range("H60").Formula = '=ROUND(((F60*(1+G60)) * $B$20/$B$21)/B60, 0)'
$ws.Range("I60").Formula = 0.815
$ws.Range("J60").Formula = '=LOG(H60*C60,2)'
